$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date value (serial 45205 -> 45206)
# for every data row from row 2 through row 122.
$startRow = 2
$endRow = 122

for ($row = $startRow; $row -le $endRow; $row++) {
    $ws.Cells.Item($row, 3).Value = 45206
}
